$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll so column D is the left-most visible column,
# and select X7 (matches the new sheetView/selection in the diff) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1

# --- Row 7: logged hours shifted out one day later ---
# P7 drops from 1.5 to 1, and the freed 0.5 moves to a new day (W7).
# Q7:V7 become "logged" (green) cells like the rest of the row, and W7
# becomes the new last logged day (also green), instead of the grey
# "not yet due" styling they had before.
$ws.Range("E7").Copy()
$ws.Range("Q7:W7").PasteSpecial(-4122)
$ws.Range("P7").Value = 1
$ws.Range("W7").Value = 0.5

# --- Row 34: more hours logged on that day ---
$ws.Range("V34").Value = 1.5

# --- Row 38 ("Zoekfilters"): task re-flagged with the same status
# colour used by rows 17/18/34, and 2 hours logged against it ---
$ws.Range("C34").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("D17").Copy()
$ws.Range("D38").PasteSpecial(-4122)

$ws.Range("V34").Copy()
$ws.Range("V38").PasteSpecial(-4122)
$ws.Range("V38").Value = 2

# --- Restore the final selection shown in the diff ---
$ws.Range("X7").Select()
